$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# --- Populate new rows 2-5 on the "Assets" sheet -----------------------
# The shared-strings table is built in the exact order the cells below are
# first written to, so new strings come out at the same indices as the
# target workbook (45..52).

# 1) "Name" column values (A) - first occurrences create the shared strings
$ws.Range("A5").Value = "EmailCredentials"
$ws.Range("A2").Value = "GoogleFormLink"
$ws.Range("A3").Value = "LetterSubject"
$ws.Range("A4").Value = "LetterText"

# 2) "Description" column values (D)
$ws.Range("D5").Value = "This is credentials for email in which we receive registration letters from users, and from which we send responses and new data about trips."
$ws.Range("D2").Value = "This is link to google form that any user of trip-planner need to fulfill to provide"
$ws.Range("D4").Value = "Text that will be send for each user who would like to start registration in trip-planner"

# 3) "Asset" column value (C) - shared by all four rows
$ws.Range("C2").Value = "Trip-Planner"

# 4) Duplicate the Name values into column B (reuses existing shared strings)
$ws.Range("B5").Value = "EmailCredentials"
$ws.Range("B2").Value = "GoogleFormLink"
$ws.Range("B3").Value = "LetterSubject"
$ws.Range("B4").Value = "LetterText"

# 5) Fill in the remaining "Asset" cells (reuse the "Trip-Planner" string)
$ws.Range("C3").Value = "Trip-Planner"
$ws.Range("C4").Value = "Trip-Planner"
$ws.Range("C5").Value = "Trip-Planner"

# --- Formatting ----------------------------------------------------------
# The "Asset"/"Description" cells (plus the whole EmailCredentials row) use
# a distinct, non-wrapping style separate from the sheet default.
$styledCells = @("C2", "D2", "C3", "C4", "D4", "A5", "B5", "C5", "D5")
foreach ($addr in $styledCells) {
    $ws.Range($addr).WrapText = $false
}

# --- Sheet view selection --------------------------------------------------
$ws.Activate()
$ws.Range("C10").Select()

# --- Page setup --------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
